$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the Lease Reference value used in B2:B8 from "OXF-002" to "OXF-TOPM2"
$ws.Range("B2:B8").Value = "OXF-TOPM2"

# 2. Fill column L (Currency = EUR) for rows 3 through 60, matching existing L2
$ws.Range("L3:L60").Value = "EUR"

# 3. Update the selection to B8
$ws.Range("B8").Select()
